# The student-ID run "STUDENT ID-AF04017789" needs to become two runs with
# the same character formatting: "STUDENT ID-AF0401" followed by "789"
# (net text also drops one duplicated "7").
$d = $word.ActiveDocument

$search = $d.Content
$found = $search.Find.Execute("STUDENT ID-AF04017789", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    # Re-materialize the hit as a plain Range (by offsets) so InsertXML
    # replaces the located text rather than being treated as an insertion
    # point.
    $target = $d.Range($search.Start, $search.End)

    $xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData>' +
           '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
           '<w:body><w:p>' +
           '<w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="56"/><w:szCs w:val="72"/></w:rPr><w:t>STUDENT ID-AF0401</w:t></w:r>' +
           '<w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="56"/><w:szCs w:val="72"/></w:rPr><w:t>789</w:t></w:r>' +
           '</w:p></w:body></w:document>' +
           '</pkg:xmlData></pkg:part></pkg:package>'

    $target.InsertXML($xml)
}
